$wb = $excel.ActiveWorkbook

# --- Part 1: Insert a new row 2 into "总计" (summary) sheet for 2022-Q3 ---
$summary = $wb.Worksheets.Item(1)

# Insert a blank row at position 2, shifting existing quarterly rows down by one.
$summary.Rows(2).Insert()

# The inserted row inherits header-row formatting; strip it back to plain, then
# restore the "index column" style (A2) by copying it from A3 (a row that still
# carries the original style).
$summary.Range("A2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# New top row: 2022-Q3 data.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 29
$summary.Range("D2").Value = 3.7

# Column A is a fixed 0-based row index (unrelated to the quarter shown in B);
# re-stamp it for every row now that one more row exists.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6

# --- Part 2: Insert new "2022-Q3" sheet with Q3 fund-holding data, before old "2022-Q2" ---
$refSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($refSheet)
$newSheet.Name = "2022-Q3"
$oldQ2 = $wb.Worksheets.Item(3)

# Header row (B1:H1)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows 2..30
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'011790"
$newSheet.Range("C2").Value = "建信创新驱动混合"
$newSheet.Range("D2").Value = "'18.42"
$newSheet.Range("E2").Value = "'89.01"
$newSheet.Range("F2").Value = "'5.16"
$newSheet.Range("G2").Value = "'0.9505"
$newSheet.Range("H2").Value = 2
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'000308"
$newSheet.Range("C3").Value = "建信创新中国混合"
$newSheet.Range("D3").Value = "'10.23"
$newSheet.Range("E3").Value = "'86.36"
$newSheet.Range("F3").Value = "'4.92"
$newSheet.Range("G3").Value = "'0.5033"
$newSheet.Range("H3").Value = 3
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'001070"
$newSheet.Range("C4").Value = "建信信息产业股票A"
$newSheet.Range("D4").Value = "'8.92"
$newSheet.Range("E4").Value = "'91.60"
$newSheet.Range("F4").Value = "'5.27"
$newSheet.Range("G4").Value = "'0.4701"
$newSheet.Range("H4").Value = 4
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'006314"
$newSheet.Range("C5").Value = "中融策略优选混合A"
$newSheet.Range("D5").Value = "'9.31"
$newSheet.Range("E5").Value = "'93.90"
$newSheet.Range("F5").Value = "'3.30"
$newSheet.Range("G5").Value = "'0.3072"
$newSheet.Range("H5").Value = 7
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'012519"
$newSheet.Range("C6").Value = "大成核心趋势混合A"
$newSheet.Range("D6").Value = "'9.97"
$newSheet.Range("E6").Value = "'88.72"
$newSheet.Range("F6").Value = "'2.80"
$newSheet.Range("G6").Value = "'0.2792"
$newSheet.Range("H6").Value = 10
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'008962"
$newSheet.Range("C7").Value = "建信科技创新混合A"
$newSheet.Range("D7").Value = "'3.69"
$newSheet.Range("E7").Value = "'91.48"
$newSheet.Range("F7").Value = "'5.69"
$newSheet.Range("G7").Value = "'0.2100"
$newSheet.Range("H7").Value = 2
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "'013561"
$newSheet.Range("C8").Value = "中融匠心优选混合A"
$newSheet.Range("D8").Value = "'5.24"
$newSheet.Range("E8").Value = "'92.24"
$newSheet.Range("F8").Value = "'3.54"
$newSheet.Range("G8").Value = "'0.1855"
$newSheet.Range("H8").Value = 8
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "'010532"
$newSheet.Range("C9").Value = "广发恒信一年持有期混合A"
$newSheet.Range("D9").Value = "'28.73"
$newSheet.Range("E9").Value = "'20.09"
$newSheet.Range("F9").Value = "'0.60"
$newSheet.Range("G9").Value = "'0.1724"
$newSheet.Range("H9").Value = 5
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "'003145"
$newSheet.Range("C10").Value = "中融竞争优势股票"
$newSheet.Range("D10").Value = "'2.86"
$newSheet.Range("E10").Value = "'94.25"
$newSheet.Range("F10").Value = "'3.18"
$newSheet.Range("G10").Value = "'0.0909"
$newSheet.Range("H10").Value = 8
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "'012520"
$newSheet.Range("C11").Value = "大成核心趋势混合C"
$newSheet.Range("D11").Value = "'2.45"
$newSheet.Range("E11").Value = "'88.72"
$newSheet.Range("F11").Value = "'2.80"
$newSheet.Range("G11").Value = "'0.0686"
$newSheet.Range("H11").Value = 10
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "'014653"
$newSheet.Range("C12").Value = "建信卓越成长一年持有混合A"
$newSheet.Range("D12").Value = "'1.97"
$newSheet.Range("E12").Value = "'91.28"
$newSheet.Range("F12").Value = "'3.38"
$newSheet.Range("G12").Value = "'0.0666"
$newSheet.Range("H12").Value = 5
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "'009135"
$newSheet.Range("C13").Value = "广发恒隆一年持有期混合A"
$newSheet.Range("D13").Value = "'5.77"
$newSheet.Range("E13").Value = "'24.66"
$newSheet.Range("F13").Value = "'0.89"
$newSheet.Range("G13").Value = "'0.0514"
$newSheet.Range("H13").Value = 6
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "'011192"
$newSheet.Range("C14").Value = "广发恒荣三个月持有期混合A"
$newSheet.Range("D14").Value = "'1.42"
$newSheet.Range("E14").Value = "'35.91"
$newSheet.Range("F14").Value = "'3.04"
$newSheet.Range("G14").Value = "'0.0432"
$newSheet.Range("H14").Value = 2
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "'009956"
$newSheet.Range("C15").Value = "广发恒誉混合A"
$newSheet.Range("D15").Value = "'3.73"
$newSheet.Range("E15").Value = "'25.82"
$newSheet.Range("F15").Value = "'1.02"
$newSheet.Range("G15").Value = "'0.0380"
$newSheet.Range("H15").Value = 9
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "'010533"
$newSheet.Range("C16").Value = "广发恒信一年持有期混合C"
$newSheet.Range("D16").Value = "'6.05"
$newSheet.Range("E16").Value = "'20.09"
$newSheet.Range("F16").Value = "'0.60"
$newSheet.Range("G16").Value = "'0.0363"
$newSheet.Range("H16").Value = 5
$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "'001189"
$newSheet.Range("C17").Value = "广发聚宝混合A"
$newSheet.Range("D17").Value = "'4.25"
$newSheet.Range("E17").Value = "'24.91"
$newSheet.Range("F17").Value = "'0.85"
$newSheet.Range("G17").Value = "'0.0361"
$newSheet.Range("H17").Value = 5
$newSheet.Range("A18").Value = 16
$newSheet.Range("B18").Value = "'006315"
$newSheet.Range("C18").Value = "中融策略优选混合C"
$newSheet.Range("D18").Value = "'1.04"
$newSheet.Range("E18").Value = "'93.90"
$newSheet.Range("F18").Value = "'3.30"
$newSheet.Range("G18").Value = "'0.0343"
$newSheet.Range("H18").Value = 7
$newSheet.Range("A19").Value = 17
$newSheet.Range("B19").Value = "'010987"
$newSheet.Range("C19").Value = "中融鑫锐研究精选一年持有期混合A"
$newSheet.Range("D19").Value = "'1.02"
$newSheet.Range("E19").Value = "'92.17"
$newSheet.Range("F19").Value = "'3.33"
$newSheet.Range("G19").Value = "'0.0340"
$newSheet.Range("H19").Value = 5
$newSheet.Range("A20").Value = 18
$newSheet.Range("B20").Value = "'009136"
$newSheet.Range("C20").Value = "广发恒隆一年持有期混合C"
$newSheet.Range("D20").Value = "'3.30"
$newSheet.Range("E20").Value = "'24.66"
$newSheet.Range("F20").Value = "'0.89"
$newSheet.Range("G20").Value = "'0.0294"
$newSheet.Range("H20").Value = 6
$newSheet.Range("A21").Value = 19
$newSheet.Range("B21").Value = "'007848"
$newSheet.Range("C21").Value = "广发聚宝混合C"
$newSheet.Range("D21").Value = "'3.33"
$newSheet.Range("E21").Value = "'24.91"
$newSheet.Range("F21").Value = "'0.85"
$newSheet.Range("G21").Value = "'0.0283"
$newSheet.Range("H21").Value = 5
$newSheet.Range("A22").Value = 20
$newSheet.Range("B22").Value = "'009937"
$newSheet.Range("C22").Value = "东方欣益一年持有期偏债混合A"
$newSheet.Range("D22").Value = "'2.22"
$newSheet.Range("E22").Value = "'20.97"
$newSheet.Range("F22").Value = "'0.72"
$newSheet.Range("G22").Value = "'0.0160"
$newSheet.Range("H22").Value = 7
$newSheet.Range("A23").Value = 21
$newSheet.Range("B23").Value = "'008963"
$newSheet.Range("C23").Value = "建信科技创新混合C"
$newSheet.Range("D23").Value = "'0.23"
$newSheet.Range("E23").Value = "'91.48"
$newSheet.Range("F23").Value = "'5.69"
$newSheet.Range("G23").Value = "'0.0131"
$newSheet.Range("H23").Value = 2
$newSheet.Range("A24").Value = 22
$newSheet.Range("B24").Value = "'010988"
$newSheet.Range("C24").Value = "中融鑫锐研究精选一年持有期混合C"
$newSheet.Range("D24").Value = "'0.35"
$newSheet.Range("E24").Value = "'92.17"
$newSheet.Range("F24").Value = "'3.33"
$newSheet.Range("G24").Value = "'0.0117"
$newSheet.Range("H24").Value = 5
$newSheet.Range("A25").Value = 23
$newSheet.Range("B25").Value = "'013562"
$newSheet.Range("C25").Value = "中融匠心优选混合C"
$newSheet.Range("D25").Value = "'0.33"
$newSheet.Range("E25").Value = "'92.24"
$newSheet.Range("F25").Value = "'3.54"
$newSheet.Range("G25").Value = "'0.0117"
$newSheet.Range("H25").Value = 8
$newSheet.Range("A26").Value = 24
$newSheet.Range("B26").Value = "'014654"
$newSheet.Range("C26").Value = "建信卓越成长一年持有混合C"
$newSheet.Range("D26").Value = "'0.23"
$newSheet.Range("E26").Value = "'91.28"
$newSheet.Range("F26").Value = "'3.38"
$newSheet.Range("G26").Value = "'0.0078"
$newSheet.Range("H26").Value = 5
$newSheet.Range("A27").Value = 25
$newSheet.Range("B27").Value = "'014863"
$newSheet.Range("C27").Value = "建信信息产业股票C"
$newSheet.Range("D27").Value = "'0.06"
$newSheet.Range("E27").Value = "'91.60"
$newSheet.Range("F27").Value = "'5.27"
$newSheet.Range("G27").Value = "'0.0032"
$newSheet.Range("H27").Value = 4
$newSheet.Range("A28").Value = 26
$newSheet.Range("B28").Value = "'009938"
$newSheet.Range("C28").Value = "东方欣益一年持有期偏债混合C"
$newSheet.Range("D28").Value = "'0.37"
$newSheet.Range("E28").Value = "'20.97"
$newSheet.Range("F28").Value = "'0.72"
$newSheet.Range("G28").Value = "'0.0027"
$newSheet.Range("H28").Value = 7
$newSheet.Range("A29").Value = 27
$newSheet.Range("B29").Value = "'011193"
$newSheet.Range("C29").Value = "广发恒荣三个月持有期混合C"
$newSheet.Range("D29").Value = "'0.07"
$newSheet.Range("E29").Value = "'35.91"
$newSheet.Range("F29").Value = "'3.04"
$newSheet.Range("G29").Value = "'0.0021"
$newSheet.Range("H29").Value = 2
$newSheet.Range("A30").Value = 28
$newSheet.Range("B30").Value = "'009957"
$newSheet.Range("C30").Value = "广发恒誉混合C"
$newSheet.Range("D30").Value = "'0.05"
$newSheet.Range("E30").Value = "'25.82"
$newSheet.Range("F30").Value = "'1.02"
$newSheet.Range("G30").Value = "'0.0005"
$newSheet.Range("H30").Value = 9

# Copy header style (B1:H1) and A-column style (A2:A30) from the old "2022-Q2" sheet
$oldQ2.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$oldQ2.Range("A2").Copy()
$newSheet.Range("A2:A30").PasteSpecial(-4122)
